$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Arun (row 4) was missing his e-mail address -- fill it in and link it as a
# mailto: hyperlink, matching every other row in the roster.
$ws.Range("B4").Value = "aruncyclopse007@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:aruncyclopse007@gmail.com")

# Hyperlinks.Add re-derives the cell's style (losing the border/fill that the
# existing "Hyperlink" cell style already carries for column B). Re-apply the
# formatting used by the other e-mail cells so B4 matches them exactly.
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)  # xlPasteFormats

# The last recorded selection moved from C21 to C22.
[void]$ws.Range("C22").Select()
